# Refresh the cryptocurrency price/volume snapshot (GitHub Actions scrape).
#
# Price/Volume(1h) cells are stored as plain text in this workbook, even
# when a value happens to look like a number (e.g. "0.4058"). Assigning
# such a string straight to Range.Value makes Excel auto-convert it to a
# real number, which would change the cell stored type. To keep it as
# text we briefly switch the cell to a text number format before the
# assignment, then restore the default "Normal" style so the cell is left
# exactly as it was otherwise (no stray formatting).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($CellRef, $Value)
    $cell = $ws.Range($CellRef)
    if ($Value -match '^[+-]?\d+(\.\d+)?$') {
        # Numeric-looking text: force a text format so it round-trips as a string.
        $cell.NumberFormat = "@"
        $cell.Value = $Value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $Value
    }
}

$updates = @(
    ,@('D2', '29.513.38')
    ,@('E2', '  +0.73%  ')
    ,@('D3', '1.913.73')
    ,@('E3', '  -0.01%  ')
    ,@('E4', '  +0.48%  ')
    ,@('D5', '325.02')
    ,@('E5', '  +0.49%  ')
    ,@('E6', '  +0.46%  ')
    ,@('E7', '  +2.17%  ')
    ,@('D8', '0.4058')
    ,@('E8', '  -0.38%  ')
    ,@('D9', '0.08138')
    ,@('E9', '  +1.45%  ')
    ,@('D10', '1.011')
    ,@('E10', '  +0.57%  ')
    ,@('D11', '23.39')
    ,@('E11', '  +4.17%  ')
    ,@('D12', '1.931.95')
    ,@('E12', '  -0.64%  ')
    ,@('D13', '5.987')
    ,@('E13', '  +1.72%  ')
    ,@('D14', '7.122')
    ,@('E14', '  -0.12%  ')
    ,@('D15', '90.17')
    ,@('E15', '  +0.58%  ')
    ,@('E16', '  +0.44%  ')
    ,@('D17', '0.06760')
    ,@('E17', '  +1.82%  ')
    ,@('D18', '0.00001039')
    ,@('E18', '  +0.99%  ')
    ,@('D19', '17.70')
    ,@('E19', '  +0.19%  ')
    ,@('E20', '  +0.54%  ')
    ,@('D21', '29.526.08')
    ,@('E21', '  +0.71%  ')
    ,@('D22', '5.627')
    ,@('E22', '  +2.21%  ')
    ,@('D23', '11.75')
    ,@('E23', '  +2.85%  ')
    ,@('E24', '  -1.04%  ')
    ,@('D25', '2.136.05')
    ,@('E25', '  -1.13%  ')
    ,@('D26', '155.72')
    ,@('E26', '  +0.80%  ')
    ,@('D27', '6.399')
    ,@('E27', '  +5.79%  ')
    ,@('D28', '20.02')
    ,@('E28', '  +1.14%  ')
    ,@('D29', '2.102')
    ,@('E29', '  -0.59%  ')
    ,@('D30', '119.83')
    ,@('E30', '  +1.92%  ')
    ,@('D31', '1.026')
    ,@('E31', '  -4.22%  ')
    ,@('D32', '0.09529')
    ,@('E32', '  -0.06%  ')
    ,@('D33', '5.513')
    ,@('E33', '  +2.24%  ')
    ,@('D34', '3.562')
    ,@('E34', '  +0.18%  ')
    ,@('D35', '1.385')
    ,@('E35', '  -2.91%  ')
    ,@('D36', '0.02269')
    ,@('E36', '  +0.85%  ')
    ,@('D37', '0.06090')
    ,@('E37', '  +0.20%  ')
    ,@('D38', '1.174')
    ,@('E38', '  +0.12%  ')
    ,@('D39', '0.5940')
    ,@('E39', '  +1.18%  ')
    ,@('D40', '7.950')
    ,@('E40', '  -3.82%  ')
    ,@('D41', '10.67')
    ,@('E41', '  +5.56%  ')
    ,@('E42', '  +0.91%  ')
    ,@('E43', '  +0.91%  ')
    ,@('D44', '2.408')
    ,@('E44', '  -5.09%  ')
    ,@('B45', 'Cronos')
    ,@('C45', 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro')
    ,@('D45', '0.07640')
    ,@('E45', '  -2.76%  ')
    ,@('B46', 'EnergySwap')
    ,@('C46', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens')
    ,@('D46', '12.48')
    ,@('E46', '  +3.31%  ')
    ,@('D47', '0.5572')
    ,@('E47', '  +0.76%  ')
    ,@('D48', '1.938')
    ,@('E48', '  +0.62%  ')
    ,@('E49', '  +3.21%  ')
    ,@('D50', '72.47')
    ,@('E50', '  +1.69%  ')
    ,@('D51', '1.053')
    ,@('E51', '  +2.12%  ')
)

foreach ($u in $updates) {
    Set-TextValue $u[0] $u[1]
}
